$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.753.07"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "2.490.82"
$ws.Range("E3").Value = "  +1.76%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.24%  "

$ws.Range("E6").Value = "  +3.08%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.25%  "

$ws.Range("D9").Value = "2.500.82"
$ws.Range("E9").Value = "  +1.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0991"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.15%  "

$ws.Range("E11").Value = "  -2.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").Value = "2.933.98"
$ws.Range("E14").Value = "  +1.78%  "

$ws.Range("D15").Value = "58.586.70"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("E17").Value = "  +1.84%  "

$ws.Range("D18").Value = "2.487.85"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("E20").Value = "  +2.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.89%  "

$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.92%  "

$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.992"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.80%  "

$ws.Range("E30").Value = "  +3.56%  "

$ws.Range("E31").Value = "  +3.86%  "

$ws.Range("E32").Value = "  +1.56%  "

$ws.Range("E33").Value = "  +1.54%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("E36").Value = "  +1.78%  "

$ws.Range("E37").Value = "  -3.02%  "

$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("E39").Value = "  +3.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.46%  "

$ws.Range("E42").Value = "  +2.60%  "

$ws.Range("E43").Value = "  +2.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "273.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.12%  "

$ws.Range("E47").Value = "  +2.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0509"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.22%  "

$ws.Range("E49").Value = "  +3.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
